$d = $word.ActiveDocument

# 1) Fix the mojibake character in the table title: "ß" -> "Ã"
$d.Content.Find.Execute(
    "Table 2. Value of the coefficients (ß) from the linear models",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Table 2. Value of the coefficients (Ã) from the linear models",
    2
) | Out-Null

# 2) Widen the 4th grid column (B_F) from 674 dxa (33.7pt) to 1007 dxa (50.35pt)
$tbl = $d.Tables.Item(1)
$tbl.Columns.Item(4).Width = 50.35

# 3) Swap the B_F / B_R values for the HAKENRTN row (row 2 of the table)
$row1 = $tbl.Rows.Item(2)
$bf1 = $row1.Cells.Item(4).Range.Text
$br1 = $row1.Cells.Item(5).Range.Text
$row1.Cells.Item(4).Range.Text = $br1
$row1.Cells.Item(5).Range.Text = $bf1

# 4) Swap the B_F / B_R values for the HAKESOTH row (row 3 of the table)
$row2 = $tbl.Rows.Item(3)
$bf2 = $row2.Cells.Item(4).Range.Text
$br2 = $row2.Cells.Item(5).Range.Text
$row2.Cells.Item(4).Range.Text = $br2
$row2.Cells.Item(5).Range.Text = $bf2
